# sign in positive scenario
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("login")

$ws.Range("B2").Value = "We_team3"
$ws.Range("A2").Value = "Ninja_Galaxy"

$ws.Range("C5").Select()
